# VS extension updated to .NET MAUI GA
#
# readme.docx: drop the "(RC3)" qualifier from the three template bullet
# points, and flip the "target .NET MAUI RC3 (...)" note to GA / Preview 1.1.
#
# Each replacement below is scoped to a Range built from the exact
# character offsets of the target substring (found via a plain string
# search inside the owning paragraph) rather than a document-wide
# Find/Replace. That keeps the edit from bleeding into neighbouring runs
# that must stay exactly as they are (e.g. the hyperlink run at the end
# of the first bullet, or the underline run wrapping "RC3 (VS2022 ...)").

function Replace-Substring($range, [string]$needle, [string]$replacement) {
    $paraText = $range.Text
    $idx = $paraText.IndexOf($needle)
    if ($idx -lt 0) {
        throw "Substring not found: $needle"
    }
    $start = $range.Start + $idx
    $end = $start + $needle.Length
    $target = $word.ActiveDocument.Range($start, $end)
    $target.Text = $replacement
}

$d = $word.ActiveDocument

# 1) ".NET MAUI App (RC3)" -> ".NET MAUI App"
Replace-Substring $d.Paragraphs(4).Range ".NET MAUI App (RC3)" ".NET MAUI App"

# 2) ".NET MAUI App (C#) (RC3)" -> ".NET MAUI App (C#)"
Replace-Substring $d.Paragraphs(5).Range ".NET MAUI App (C#) (RC3)" ".NET MAUI App (C#)"

# 3) ".NET MAUI Class Library (RC3)" -> ".NET MAUI Class Library"
Replace-Substring $d.Paragraphs(6).Range ".NET MAUI Class Library (RC3)" ".NET MAUI Class Library"

# 4) ".NET MAUI RC3 (VS2022 17.3 Preview 1.0 or later)"
#    -> ".NET MAUI GA (VS2022 17.3 Preview 1.1 or later)"
#    (this whole span carries single-underline formatting in the source;
#    keeping the Range confined to it preserves that formatting)
Replace-Substring $d.Paragraphs(219).Range `
    ".NET MAUI RC3 (VS2022 17.3 Preview 1.0 or later)" `
    ".NET MAUI GA (VS2022 17.3 Preview 1.1 or later)"
